$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# --- Sheet "choices" (sheet2): reorder / replace the holiday choice list ---
$choices.Range("B19").Value = "kwanzaa"
$choices.Range("D19").Value = "Kwanzaa"

$choices.Range("B20").Value = "christmas"
$choices.Range("D20").Value = "Christmas"

$choices.Range("B21").Value = "hannukah"
$choices.Range("D21").Value = "Hannukah"

$choices.Range("D22").Value = "Diwali"
$choices.Range("B22").Value = "diwali"

# --- Sheet "survey" (sheet1): add "//" comment markers to if/end-if and
# begin/end screen rows, used for commenting-out lines during form development ---
$survey.Range("B21").Value = "//if"
$survey.Range("B22").Value = "//"
$survey.Range("B23").Value = "// end if"

$survey.Range("B27").Value = "//begin screen"
$survey.Range("B28").Value = "//"
$survey.Range("B29").Value = "//"
$survey.Range("B30").Value = "//"
$survey.Range("B31").Value = "//"
$survey.Range("B32").Value = "//end screen"

# --- View / selection changes ---

# choices sheet selection -> E20
[void]$choices.Activate()
[void]$choices.Range("E20").Select()

# queries sheet keeps its own selection (C7), but loses tabSelected
$queries = $wb.Worksheets.Item("queries")
[void]$queries.Activate()
[void]$queries.Range("C7").Select()

# survey sheet becomes the active tab, with selection -> B33
[void]$survey.Activate()
[void]$survey.Range("B33").Select()
